$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.398.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.38%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.592.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.01%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.76%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "568.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.07%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.09%  "
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.22%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.681"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.61%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "63.93"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +15.39%  "
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.40%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.84%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.45%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.169.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.01%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.593.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.87%  "
# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.21"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.04%  "
# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.126"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.33%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.219.48"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.18%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.22"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.18%  "
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.60%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "405.59"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.97%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.61%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.33"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.54%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.97%  "
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.39%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.12%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.83"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.26%  "
# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "744.66"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +16.84%  "
# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.40%  "
# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.85%  "
# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.57%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.90%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.19%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.62"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.64%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.422"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.46%  "
# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.21%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.54%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.19"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +31.19%  "
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.45%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.164.69"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.70%  "
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.33%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.06%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.98%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.51%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.65%  "
# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.45%  "
# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.10%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.86%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.73"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.96%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.48%  "
